# Update the cryptos worksheet with freshly scraped values (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that sometimes look like plain numbers
# (e.g. "247.79"). Force the whole price column to Text first so Excel
# doesn't silently coerce those values into numeric cells, then clear the
# number-format override again once all values are in place so the column
# is left in its original (default) style.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "36.976.81"
$ws.Range("E2").Value = "  -0.19%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.046.35"
$ws.Range("E3").Value = "  +0.14%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.41%  "

# Row 5 - BNB
$ws.Range("D5").Value = "247.79"
$ws.Range("E5").Value = "  +0.26%  "

# Row 6 - XRP
$ws.Range("D6").Value = "0.660"
$ws.Range("E6").Value = "  +1.49%  "

# Row 7 - now USDC (was Solana)
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.10%  "

# Row 8 - now Solana (was USDC)
$ws.Range("B8").Value = "Solana"
$ws.Range("C8").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D8").Value = "57.39"
$ws.Range("E8").Value = "  +5.04%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.379"
$ws.Range("E9").Value = "  +1.47%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.0775"
$ws.Range("E10").Value = "  +2.42%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +1.89%  "

# Row 12 - Chainlink
$ws.Range("D12").Value = "15.70"
$ws.Range("E12").Value = "  +5.47%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "2.346.26"
$ws.Range("E13").Value = "  +0.30%  "

# Row 14 - Polygon
$ws.Range("D14").Value = "0.797"
$ws.Range("E14").Value = "  -1.42%  "

# Row 15 - Polkadot
$ws.Range("D15").Value = "5.53"
$ws.Range("E15").Value = "  +7.45%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "2.045.34"
$ws.Range("E16").Value = "  +0.15%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "37.042.10"

# Row 18 - Avalanche
$ws.Range("D18").Value = "16.46"
$ws.Range("E18").Value = "  +17.09%  "

# Row 19 - Litecoin
$ws.Range("D19").Value = "74.35"
$ws.Range("E19").Value = "  +3.91%  "

# Row 20 - ShibaInu (price uses a subscript-3 digit-grouping glyph, U+2083)
$shibaSub3 = [char]0x2083
$ws.Range("D20").Value = "0.0{0}0899" -f $shibaSub3
$ws.Range("E20").Value = "  +1.38%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "5.31"
$ws.Range("E21").Value = "  +2.24%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "235.11"
$ws.Range("E22").Value = "  -0.18%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  -0.07%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -1.33%  "

# Row 25 - PancakeSwap
$ws.Range("D25").Value = "2.17"
$ws.Range("E25").Value = "  +10.45%  "

# Row 26 - Monero
$ws.Range("D26").Value = "167.50"
$ws.Range("E26").Value = "  -0.73%  "

# Row 27 - Cosmos
$ws.Range("D27").Value = "9.07"
$ws.Range("E27").Value = "  +1.28%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "19.65"
$ws.Range("E28").Value = "  -1.20%  "

# Row 29 - Stellar
$ws.Range("D29").Value = "0.123"
$ws.Range("E29").Value = "  +1.49%  "

# Row 30 - ImmutableX
$ws.Range("D30").Value = "1.12"
$ws.Range("E30").Value = "  +6.71%  "

# Row 31 - Filecoin
$ws.Range("D31").Value = "4.66"
$ws.Range("E31").Value = "  +3.55%  "

# Row 32 - Hedera
$ws.Range("D32").Value = "0.0609"
$ws.Range("E32").Value = "  -0.92%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("D33").Value = "4.44"
$ws.Range("E33").Value = "  +3.42%  "

# Row 34 - Kaspa
$ws.Range("D34").Value = "0.0879"
$ws.Range("E34").Value = "  +0.98%  "

# Row 35 - BinanceUSD
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.19%  "

# Row 36 - LidoDAOToken
$ws.Range("D36").Value = "2.20"
$ws.Range("E36").Value = "  -1.68%  "

# Row 37 - WEMIXToken
$ws.Range("E37").Value = "  -1.96%  "

# Row 38 - Cronos
$ws.Range("D38").Value = "0.106"
$ws.Range("E38").Value = "  +2.84%  "

# Row 39 - TrustWalletToken
$ws.Range("E39").Value = "  +0.28%  "

# Row 40 - HuobiToken
$ws.Range("E40").Value = "  +14.93%  "

# Row 41 - THORChain
$ws.Range("E41").Value = "  +25.07%  "

# Row 42 - VeChain
$ws.Range("D42").Value = "0.0219"
$ws.Range("E42").Value = "  -0.89%  "

# Row 43 - InjectiveProtocol
$ws.Range("D43").Value = "17.20"
$ws.Range("E43").Value = "  -5.20%  "

# Row 44 - ARBITRUM
$ws.Range("D44").Value = "1.11"
$ws.Range("E44").Value = "  -1.12%  "

# Row 45 - Aave
$ws.Range("D45").Value = "95.12"
$ws.Range("E45").Value = "  +0.21%  "

# Row 46 - RenderToken
$ws.Range("D46").Value = "2.42"
$ws.Range("E46").Value = "  +4.17%  "

# Row 47 - Maker
$ws.Range("D47").Value = "1.274.60"
$ws.Range("E47").Value = "  -0.89%  "

# Row 48 - MXToken
$ws.Range("E48").Value = "  -1.57%  "

# Row 49 - RocketPoolETH
$ws.Range("D49").Value = "2.231.20"
$ws.Range("E49").Value = "  +0.97%  "

# Row 50 - FraxShare
$ws.Range("D50").Value = "6.64"
$ws.Range("E50").Value = "  -0.87%  "

# Row 51 - FTXToken
$ws.Range("D51").Value = "3.54"
$ws.Range("E51").Value = "  -12.70%  "

# Restore the default (General) number format on the price column now that
# every value has been written as text, so the saved style matches the
# original workbook (no lingering text-format override).
$ws.Range("D2:D51").ClearFormats()
